$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.882.06"
$ws.Range("E2").Value = "  +3.54%  "

$ws.Range("D3").Value = "3.131.37"
$ws.Range("E3").Value = "  +2.37%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "590.98"
$ws.Range("E5").Value = "  +2.73%  "

$ws.Range("D6").Value = "146.69"
$ws.Range("E6").Value = "  +2.94%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "3.123.51"
$ws.Range("E8").Value = "  +2.28%  "

$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("E10").Value = "  +19.38%  "

$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  +4.99%  "

$ws.Range("D12").Value = "0.468"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  +7.72%  "

$ws.Range("D14").Value = "36.14"
$ws.Range("E14").Value = "  +3.67%  "

$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").Value = "3.649.33"
$ws.Range("E16").Value = "  +2.37%  "

$ws.Range("D17").Value = "7.20"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "63.796.52"
$ws.Range("E18").Value = "  +3.34%  "

$ws.Range("D19").Value = "3.131.95"
$ws.Range("E19").Value = "  +2.30%  "

$ws.Range("D20").Value = "465.76"
$ws.Range("E20").Value = "  +3.49%  "

$ws.Range("D21").Value = "14.28"
$ws.Range("E21").Value = "  +2.65%  "

$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  +0.62%  "

$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +2.94%  "

$ws.Range("D24").Value = "13.30"
$ws.Range("E24").Value = "  -1.55%  "

$ws.Range("D25").Value = "82.52"
$ws.Range("E25").Value = "  +0.85%  "

$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("E27").Value = "  +8.36%  "

$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  +3.27%  "

$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "6.84"
$ws.Range("E31").Value = "  +2.93%  "

$ws.Range("E32").Value = "  +2.00%  "

$ws.Range("E33").Value = "  +1.50%  "

$ws.Range("D34").Value = "0.0₃0872"
$ws.Range("E34").Value = "  +8.52%  "

$ws.Range("E35").Value = "  +10.91%  "

$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +2.31%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "3.44"
$ws.Range("E37").Value = "  +16.09%  "

$ws.Range("D38").Value = "6.15"
$ws.Range("E38").Value = "  +1.86%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "50.86"
$ws.Range("E39").Value = "  +1.60%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "453.38"
$ws.Range("E40").Value = "  +7.81%  "

$ws.Range("E41").Value = "  -0.83%  "

$ws.Range("D42").Value = "0.0372"
$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").Value = "2.922.13"
$ws.Range("E43").Value = "  +5.47%  "

$ws.Range("E44").Value = "  +4.54%  "

$ws.Range("E45").Value = "  +2.63%  "

$ws.Range("E46").Value = "  +3.66%  "

$ws.Range("D47").Value = "128.41"
$ws.Range("E47").Value = "  +4.22%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").Value = "34.58"
$ws.Range("E49").Value = "  -4.94%  "

$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").Value = "24.66"
$ws.Range("E51").Value = "  +3.22%  "
